$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the bold / thin-bordered / center+top-aligned style once on B1, then
# clone it onto A2 via a single copy -> paste-special(formats) so the style
# table ends up with exactly one extra cellXf (shared by both cells) instead
# of accumulating unused intermediate style records.
$c1 = $ws.Range("B1")
$c1.Font.Bold = $true
$c1.Borders.LineStyle = 1
$c1.HorizontalAlignment = -4108  # xlCenter
$c1.VerticalAlignment = -4160    # xlTop

$c1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Cell values: B1 and A2 are numeric zeros, B2 is the text label.
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"
